$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels for start year 2015 -> 2016
$ws.Range("B1").Value = "Pop2016"
$ws.Range("C1").Value = "UrbanRatio2016"

# Update data values in row 2
$ws.Range("B2").Value = 10872000
$ws.Range("C2").Value = 0.44395000000000001
$ws.Range("R2").Value = 0.41403000000000001

# Update the active cell selection to match the author's edit
$ws.Range("E3").Select()
